$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

$ws.Cells.Item($row, 1).Value = "'-484"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "'6/18/2025"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "Conde 2319"

$ws.Cells.Item($row, 4).Value = "'13"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = "'807605771"
$ws.Cells.Item($row, 5).Style = "Normal"

$ws.Cells.Item($row, 6).Value = "GESTION TELECENTRO"

$ws.Cells.Item($row, 7).Value = "Pendiente"

$ws.Cells.Item($row, 8).Value = "'"
$ws.Cells.Item($row, 8).Style = "Normal"

$ws.Cells.Item($row, 9).Value = "'1"
$ws.Cells.Item($row, 9).Style = "Normal"

$ws.Cells.Item($row, 10).Value = "Cambio"

$ws.Cells.Item($row, 11).Value = "Fuente TLC"

$ws.Cells.Item($row, 12).Value = "Pasante"

$ws.Cells.Item($row, 13).Value = -58.467271

$ws.Cells.Item($row, 14).Value = -34.56515
